# Update countries & provincias Spain
#
# The source data for Huelva/Huesca was refreshed and the two provinces
# swapped places in the feed; "Casos activos" (column C) values also swap
# between the two rows (Huelva: 72 -> 0, Huesca: 0 -> 72), while the other
# columns (B, D, E) stay identical for both rows. The "last updated" banner
# in A1 is bumped to the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 was Huelva (37, 72, 37, 0) -> becomes Huesca (37, 0, 37, 0)
$ws.Range("A53").Value = "Huesca"
$ws.Range("C53").Value = 0

# Row 54 was Huesca (37, 0, 37, 0) -> becomes Huelva (37, 72, 37, 0)
$ws.Range("A54").Value = "Huelva"
$ws.Range("C54").Value = 72

# Bump the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 00:16"
